# Apply updated cryptocurrency price / 1h-volume data (and a couple of row swaps)
# to the active worksheet, mirroring the automated GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.102.37"
$ws.Range("E2").Value = "  +2.10%  "
$ws.Range("D3").Value = "2.513.70"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.21"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.96"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.537"
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "2.516.29"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.161"
$ws.Range("E11").Value = "  -1.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.42"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +1.62%  "
$ws.Range("D16").Value = "2.974.22"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "63.991.19"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "2.525.71"
$ws.Range("E18").Value = "  +2.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.83"
$ws.Range("E19").Value = "  -2.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.96"
$ws.Range("E20").Value = "  -0.21%  "
$ws.Range("E21").Value = "  +2.82%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "327.77"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.26"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.04"
$ws.Range("E25").Value = "  -1.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "65.37"
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "649.83"
$ws.Range("E27").Value = "  -1.35%  "
$ws.Range("E28").Value = "  +5.42%  "
$ws.Range("D29").Value = "2.642.71"
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.50"
$ws.Range("E30").Value = "  +4.02%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +1.56%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.54"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.52"
$ws.Range("E38").Value = "  +2.18%  "
$ws.Range("E39").Value = "  +0.76%  "
$ws.Range("E40").Value = "  +0.97%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "152.49"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.81"
$ws.Range("E42").Value = "  +2.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.78"
$ws.Range("E43").Value = "  +2.47%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("E45").Value = "  +6.15%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "0.0₆0301"
$ws.Range("E47").Value = "  -2.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.43"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "21.10"
$ws.Range("E50").Value = "  +2.87%  "
$ws.Range("E51").Value = "  +1.55%  "
